$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.27246105670929
$ws.Range("B1").Value = 2.926095724105835
$ws.Range("C1").Value = 5.461709499359131
$ws.Range("D1").Value = 1.86314857006073
$ws.Range("E1").Value = 1.026785016059875
